# Append a new WeightMeasurements row (WeightID 6 / 28 Jan 2018) that
# duplicates the prior week's readings (same Stone/Kg/lbs weight as row 6),
# and extend the Gain/Loss + TotalToTarget formulas down into it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 43128
$ws.Range("C7").Value = "15.11"
$ws.Range("D7").Value = 100.2
$ws.Range("E7").Value = 220.9
$ws.Range("F7").Formula = "=E7-E6"
$ws.Range("G7").Formula = "=E7-210"
